$wb = $excel.ActiveWorkbook

$wsGroup = $wb.Worksheets.Item("Group definitions")
$wsStoch = $wb.Worksheets.Item("Stochastic")
$wsVoices = $wb.Worksheets.Item("Voices")

# ---------------------------------------------------------------
# Sheet "Group definitions": update the custom Vs-law formula text
# ---------------------------------------------------------------
$wsGroup.Range("H2").Value = "2*G + G^0.25 + 3; D - D^0.6 + 10"

# ---------------------------------------------------------------
# Sheet "Stochastic": new "Analys type" column (S)
# ---------------------------------------------------------------
$wsStoch.Range("Q1").Copy() | Out-Null
$wsStoch.Range("S1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
$wsStoch.Range("S1").Value = "Analys type"
$wsStoch.Columns.Item(19).ColumnWidth = 9.9

# Bold + center header style for A1 (was bold only)
$wsStoch.Range("A1").HorizontalAlignment = -4108   # xlCenter

# ---------------------------------------------------------------
# Sheet "Voices": new helper list (MOPS / Permutations) in column C
# ---------------------------------------------------------------
$wsVoices.Range("C2").Value = "Permutations"
$wsVoices.Range("C1").Value = "MOPS"
$wsVoices.Columns.Item(3).ColumnWidth = 11.1

# ---------------------------------------------------------------
# Sheet "Stochastic": new "Brick size" column (L)
# ---------------------------------------------------------------
$wsStoch.Range("K1").Copy() | Out-Null
$wsStoch.Range("L1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
$wsStoch.Range("L1").Value = "Brick size`n[m]"
$wsStoch.Range("L2").Value = 3
$wsStoch.Columns.Item(12).ColumnWidth = 13.7

# Number of iterations changed, plus new "Random seed" value, and the
# selected "Analys type" for this run
$wsStoch.Range("O2").Value = 200
$wsStoch.Range("P2").Value = 7
$wsStoch.Range("S2").Value = "MOPS"

# ---------------------------------------------------------------
# Data validation list for the new "Analys type" cell (S2)
# ---------------------------------------------------------------
$wsStoch.Range("S2").Validation.Add(3, 1, 3, "=Voices!`$C`$1:`$C`$2") | Out-Null

# ---------------------------------------------------------------
# Selections (as left by the author when saving)
# ---------------------------------------------------------------
$wsVoices.Range("C28").Select() | Out-Null
$wsGroup.Range("A4").Select() | Out-Null
$wsStoch.Activate()
$wsStoch.Range("E4").Select() | Out-Null
